$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 0.001180555555555556
$ws.Range("K2").Value = 5381
$ws.Range("L2").Value = 0.010762
